# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header row (row 1): new labels, styled the same as the other headers.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header style (bold / centered / bordered, same as A1:AC1)
# onto the three new header cells so they reuse style index 1 rather than
# creating a near-duplicate style.
$ws.Cells.Item(1, 1).Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# -- Data rows 2-46: every player on the roster shares the team's season
# record (87 wins, 75 losses, 0 ties) for this file.
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 87  # AD - Wins
    $ws.Cells.Item($row, 31).Value = 75  # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF - Ties
}
